# "Update US model copy" — toggle which electricity source gets guaranteed
# dispatch priority in the DPbES sheet: turn "hard coal" OFF (1 -> 0) and
# turn "hard coal w CCS" ON (0 -> 1) for every forecast year column (B:AE),
# then leave the sheet's selection where the author left it (on the newly
# edited "hard coal" row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DPbES")
$ws.Activate()

# Row 2 = "hard coal": was preferenced (1) in every year, now not (0).
$ws.Range("B2:AE2").Value = 0

# Row 19 = "hard coal w CCS": was not preferenced (0), now is (1).
$ws.Range("B19:AE19").Value = 1

# Scroll/select so column D is the left-most visible column and the active
# selection sits on row 2 (B2:AE2), matching where the author left the
# cursor after making the edit.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("B2:AE2").Select()
